$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "Hep4_971_941_041122.xlsx"
$ws.Range("B34").Value = "All Hep04 data"
$ws.Range("C34").Value = 1
$ws.Range("C34").HorizontalAlignment = -4131
$ws.Range("D34").Value = "041122"
$ws.Range("E34").Value = 971
$ws.Range("E34").HorizontalAlignment = -4131
$ws.Range("F34").Value = "m4:2 FTOH"
$ws.Range("G34").Value = "Name"
$ws.Range("G34").HorizontalAlignment = -4131
$ws.Range("H34").Value = "Type"
$ws.Range("H34").HorizontalAlignment = -4131
$ws.Range("I34").Value = "Area...13"
$ws.Range("J34").Value = "Area...17"
$ws.Range("J34").HorizontalAlignment = -4131
$ws.Range("K34").Value = "CC Concentration (nM)"
$ws.Range("K34").HorizontalAlignment = -4131
$ws.Range("L34").Value = "RT...12"
$ws.Range("L34").HorizontalAlignment = -4131

$ws.Range("F10").Select() | Out-Null
